# Poloniex history workbook — "updated in 10 min"
# Row 66's buy order finalizes (Status -> DONE, Finalized date + Fee filled in),
# and a brand-new "IN PROGRESS" sell row (67) is appended right after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 66: mark the pending buy as finalized ---
$ws.Range("H66").Value = "DONE"
$ws.Range("I66").Value = 42876.604861111111
$ws.Range("J66").Value = "0.20100000 XRP (0.15%)"

# --- Row 67: new sell transaction row ---
# Clone row 66's cell formatting first (A:I and K), so the new row keeps the
# same number formats / wrap / styles as its neighbour, then overwrite values.
$ws.Range("A66:I66").Copy()
$ws.Range("A67:I67").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("K66").Copy()
$ws.Range("K67").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("A67").Value = "                5/21/2017  14:49:35 AM"
$ws.Range("B67").Value = $ws.Range("B65").Value()
$ws.Range("C67").Value = $ws.Range("C66").Value()
$ws.Range("D67").Value = 0.33700000000000002
$ws.Range("E67").Value = "         0.367  USDT"
$ws.Range("F67").Value = $ws.Range("F66").Value()
$ws.Range("G67").Value = $ws.Range("G66").Value()
$ws.Range("H67").Value = "IN PROGRESS"
$ws.Range("K67").Value = $ws.Range("K66").Value()

$ws.Rows.Item(67).RowHeight = 14.25

# Reflect the author's final selection/view state
$ws.Range("F73").Select()
